$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 7-11 are cyclically shifted up by one position: the values that were in
# row 8 move to row 7, row 9 -> row 8, row 10 -> row 9, row 11 -> row 10, and the
# original row 7 values wrap around into row 11.
# Columns B (Class) and C (First_Detection_Date) are identical across these rows,
# so only columns A and D:J need to move.
#
# We use Copy / PasteSpecial (instead of reading/writing .Value as plain strings)
# so that text-like values (e.g. "0.76", "702,633,740,690") keep their original
# text data type instead of being auto-converted to numbers, and so that no new
# cell styles get introduced.

# 1) Stash the original row 7 contents in an unused scratch area off to the side.
$ws.Range("A7").Copy()
$ws.Range("L7").PasteSpecial()
$ws.Range("D7:J7").Copy()
$ws.Range("M7:S7").PasteSpecial()

# 2) Shift rows 8,9,10,11 up into rows 7,8,9,10.
for ($r = 7; $r -le 10; $r++) {
    $srcRow = $r + 1
    $ws.Range("A$srcRow").Copy()
    $ws.Range("A$r").PasteSpecial()
    $ws.Range("D$srcRow" + ":J$srcRow").Copy()
    $ws.Range("D$r" + ":J$r").PasteSpecial()
}

# 3) Restore the stashed original row 7 values into row 11 (wrap-around).
$ws.Range("L7").Copy()
$ws.Range("A11").PasteSpecial()
$ws.Range("M7:S7").Copy()
$ws.Range("D11:J11").PasteSpecial()

# 4) Clean up the scratch area.
$ws.Range("L7:S7").ClearContents()
